# Generate Report for Handoff
# Update status text and timestamps on all three sheets to reflect the
# report having moved from "In Translation" to "Ready for handoff".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status columns (E2, F2) + HO generate date (G2) ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 00:42:54"

# --- zh-cn sheet: Status (C2), Latest Handoff Datetime (H2) ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 00:42:49"

# --- de-de sheet: Status (C2) ---
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Column width adjustments (widening to fit the longer "Ready for handoff" text) ---
# Target stored width (OOXML <col width>) is 17.2159881591797 characters; the
# COM ColumnWidth setter here quantizes to the nearest 1/6 character, so use
# the character count that lands on the nearest representable width.
$newWidth = 16.333333333333336
$wsOverview.Range("E1").ColumnWidth = $newWidth
$wsOverview.Range("F1").ColumnWidth = $newWidth
$wsZhCn.Range("C1").ColumnWidth = $newWidth
$wsDeDe.Range("C1").ColumnWidth = $newWidth
